$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Predicted_Signal (AC) and Actual_Return (AD) values to reflect
# results from the refactored (split into smaller modules) test/train pipeline.
$ws.Range("AC2").Value = 0
$ws.Range("AC3").Value = 1
$ws.Range("AD3").Value = -0
$ws.Range("AD4").Value = -0.009458618543870534
$ws.Range("AC7").Value = 1
